$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1484.1111
$ws.Range("J2").Value = 1138.8572
$ws.Range("L2").Value = 1138.8572
$ws.Range("N2").Value = -1364.8572
$ws.Range("H28").Value = 681.53125
$ws.Range("J28").Value = 1480.5555
$ws.Range("L28").Value = 1480.5555
$ws.Range("N28").Value = -2450.5555
$ws.Range("H33").Value = 353.1
$ws.Range("I33").Value = 133
$ws.Range("K33").Value = 133
$ws.Range("M33").Value = 96
$ws.Range("H111").Value = 500
$ws.Range("J111").Value = 1000
$ws.Range("L111").Value = 3000
$ws.Range("N111").Value = -9134
$ws.Range("H132").Value = 15350.16
$ws.Range("I132").Value = 15658.062
$ws.Range("J132").Value = 263
$ws.Range("K132").Value = 46974.186
$ws.Range("L132").Value = 789
$ws.Range("M132").Value = -44444.186
$ws.Range("N132").Value = -5849
$ws.Range("H135").Value = 4405.6665
$ws.Range("I135").Value = 3787.2222
$ws.Range("K135").Value = 34084.99980000001
$ws.Range("M135").Value = -31549.99980000001
$ws.Range("H137").Value = 27799
$ws.Range("I137").Value = 52332.168
$ws.Range("J137").Value = 3265.8333
$ws.Range("K137").Value = 156996.504
$ws.Range("L137").Value = 9797.499899999999
$ws.Range("M137").Value = -154446.504
$ws.Range("N137").Value = -14897.4999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9111.571
$ws.Range("I61").Value = 1305.3334
$ws.Range("J61").Value = 23162.8
$ws.Range("K61").Value = 1305.3334
$ws.Range("L61").Value = 23162.8
$ws.Range("M61").Value = -1093.3334
$ws.Range("N61").Value = -23586.8
$ws.Range("H74").Value = 279004.88
$ws.Range("I74").Value = 376055
$ws.Range("K74").Value = 376055
$ws.Range("M74").Value = -375181
$ws.Range("H77").Value = 279004.88
$ws.Range("I77").Value = 376055
$ws.Range("K77").Value = 1880275
$ws.Range("M77").Value = -1875907
$ws.Range("H122").Value = 3259.5356
$ws.Range("I122").Value = 3187.6843
$ws.Range("K122").Value = 9563.052899999999
$ws.Range("M122").Value = -7113.052899999999
$ws.Range("H136").Value = 9111.571
$ws.Range("I136").Value = 1305.3334
$ws.Range("J136").Value = 23162.8
$ws.Range("K136").Value = 3916.0002
$ws.Range("L136").Value = 69488.39999999999
$ws.Range("M136").Value = -1366.0002
$ws.Range("N136").Value = -74588.39999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14495.5
$ws.Range("I20").Value = 27035.584
$ws.Range("J20").Value = 1955.4166
$ws.Range("K20").Value = 27035.584
$ws.Range("L20").Value = 1955.4166
$ws.Range("M20").Value = -26788.584
$ws.Range("N20").Value = -2449.4166
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 68366.39999999999
$ws.Range("I132").Value = 84333
$ws.Range("K132").Value = 252999
$ws.Range("M132").Value = -250469
$ws.Range("H134").Value = 2251.2083
$ws.Range("I134").Value = 2110.318
$ws.Range("K134").Value = 6330.954000000001
$ws.Range("M134").Value = -3795.954000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 789
$ws.Range("I5").Value = 800.25
$ws.Range("K5").Value = 2400.75
$ws.Range("M5").Value = -2288.75
$ws.Range("H28").Value = 42602.6
$ws.Range("J28").Value = 9999
$ws.Range("L28").Value = 29997
$ws.Range("N28").Value = -30461
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H132").Value = 1381.9286
$ws.Range("J132").Value = 1891.8889
$ws.Range("L132").Value = 17027.0001
$ws.Range("N132").Value = -22087.0001
$ws.Range("H135").Value = 789
$ws.Range("I135").Value = 800.25
$ws.Range("K135").Value = 7202.25
$ws.Range("M135").Value = -4667.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 448.26666
$ws.Range("I2").Value = 536.0526
$ws.Range("K2").Value = 536.0526
$ws.Range("M2").Value = -423.0526
$ws.Range("H80").Value = 4922.645
$ws.Range("I80").Value = 3659.35
$ws.Range("J80").Value = 7219.5454
$ws.Range("K80").Value = 3659.35
$ws.Range("L80").Value = 7219.5454
$ws.Range("M80").Value = -2661.35
$ws.Range("N80").Value = -9215.545399999999
$ws.Range("H83").Value = 4922.645
$ws.Range("I83").Value = 3659.35
$ws.Range("J83").Value = 7219.5454
$ws.Range("K83").Value = 18296.75
$ws.Range("L83").Value = 36097.727
$ws.Range("M83").Value = -13304.75
$ws.Range("N83").Value = -46081.727
$ws.Range("H113").Value = 3237.6667
$ws.Range("I113").Value = 1603.25
$ws.Range("J113").Value = 6506.5
$ws.Range("K113").Value = 1603.25
$ws.Range("L113").Value = 6506.5
$ws.Range("M113").Value = 566.75
$ws.Range("N113").Value = -10846.5
$ws.Range("I122").Value = 2860.5833
$ws.Range("K122").Value = 8581.749899999999
$ws.Range("M122").Value = -6131.749899999999
$ws.Range("H126").Value = 1489.3334
$ws.Range("I126").Value = 1525.25
$ws.Range("K126").Value = 4575.75
$ws.Range("M126").Value = -2105.75
$ws.Range("H132").Value = 2515.0908
$ws.Range("I132").Value = 2374.2104
$ws.Range("J132").Value = 3407.3333
$ws.Range("K132").Value = 7122.6312
$ws.Range("L132").Value = 10221.9999
$ws.Range("M132").Value = -4592.6312
$ws.Range("N132").Value = -15281.9999
$ws.Range("H136").Value = 21889.666
$ws.Range("J136").Value = 21889.666
$ws.Range("L136").Value = 65668.99800000001
$ws.Range("N136").Value = -70768.99800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4943.8887
$ws.Range("I7").Value = 3785.5715
$ws.Range("K7").Value = 3785.5715
$ws.Range("M7").Value = -3673.5715
$ws.Range("H16").Value = 3500.2
$ws.Range("I16").Value = 3125.25
$ws.Range("K16").Value = 3125.25
$ws.Range("M16").Value = -2955.25
$ws.Range("H22").Value = 1310.091
$ws.Range("I22").Value = 1182.4615
$ws.Range("J22").Value = 1494.4445
$ws.Range("K22").Value = 1182.4615
$ws.Range("L22").Value = 1494.4445
$ws.Range("M22").Value = -887.4614999999999
$ws.Range("N22").Value = -2084.4445
$ws.Range("H25").Value = 46599.332
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 46599.332
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 46599.332
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -47059.332
$ws.Range("H27").Value = 1310.091
$ws.Range("I27").Value = 1182.4615
$ws.Range("J27").Value = 1494.4445
$ws.Range("K27").Value = 1182.4615
$ws.Range("L27").Value = 1494.4445
$ws.Range("M27").Value = -1075.4615
$ws.Range("N27").Value = -1708.4445
$ws.Range("H61").Value = 731.5
$ws.Range("I61").Value = 731.5
$ws.Range("K61").Value = 731.5
$ws.Range("M61").Value = -529.5
$ws.Range("H68").Value = 5285.7856
$ws.Range("I68").Value = 5199.8
$ws.Range("J68").Value = 5333.5557
$ws.Range("K68").Value = 5199.8
$ws.Range("L68").Value = 5333.5557
$ws.Range("M68").Value = -4450.8
$ws.Range("N68").Value = -6831.5557
$ws.Range("H71").Value = 5285.7856
$ws.Range("I71").Value = 5199.8
$ws.Range("J71").Value = 5333.5557
$ws.Range("K71").Value = 25999
$ws.Range("L71").Value = 26667.7785
$ws.Range("M71").Value = -22255
$ws.Range("N71").Value = -34155.7785
$ws.Range("H113").Value = 731.5
$ws.Range("I113").Value = 731.5
$ws.Range("K113").Value = 731.5
$ws.Range("M113").Value = 1438.5
$ws.Range("H126").Value = 4943.8887
$ws.Range("I126").Value = 3785.5715
$ws.Range("K126").Value = 11356.7145
$ws.Range("M126").Value = -8886.7145
$ws.Range("H132").Value = 5092.5
$ws.Range("I132").Value = 4856
$ws.Range("K132").Value = 14568
$ws.Range("M132").Value = -12038

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 177526
$ws.Range("I132").Value = 234701.33
$ws.Range("K132").Value = 704103.99
$ws.Range("M132").Value = -701573.99
